$d = $word.ActiveDocument

# --- Change 1: merge the "5.  " / "What is the difference..." runs into one run ---
$d.Content.Find.Execute(
    "5.  What is the difference between library and framework.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "5.  What is the difference between library and framework.", 2) | Out-Null

# --- Change 2: turn the lone "23. " paragraph into question 23, then add
#     the new questions 24-30 (with two blank separator paragraphs) right
#     after it, before the existing trailing blank paragraphs. ---
$d.Content.Find.Execute(
    "23. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "23. What is a components and what it contains.", 2) | Out-Null

# locate the paragraph we just renamed so we can fix its font size and
# use it as the anchor to insert the following paragraphs after.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "23. What is a components and what it contains.*") {
        $anchor = $p
        break
    }
}
$anchor.Range.Font.Size = 13

$newTexts = @(
    "24. How trhe functional and class component returns the markup.",
    "25. In which file we need to import boostrap.",
    "26. What are the JSX rules are there in React.",
    "",
    "27. How use Arraow functions in component.",
    "28. When to use normal function and when to use arrow function.",
    "29. What is the benfites of using arrow function.",
    "",
    "30. "
)

foreach ($text in $newTexts) {
    $anchor.Range.InsertParagraphAfter()
    $anchorIndex = $anchor.Range.Information(3)
    $nextPara = $anchor.Next()
    if ($text -ne "") {
        $nextPara.Range.Text = $text
    }
    $nextPara.Range.Font.Size = 13
    $anchor = $nextPara
}
